$d = $word.ActiveDocument

# Locate the run of text that needs to be split: the sentence ending in the
# misspelled "seting." (which should become "setting.") inside the
# "I would be a good fit..." paragraph.
$oldText = "with many different tools and software such as Python, Solidworks, and computer vision in a manner and depth that are unavailable in an academic seting. "

$full = $d.Content.Text
$startIdx = $full.IndexOf($oldText)
$endIdx = $startIdx + $oldText.Length

# Remove the whole original run's text so the following inserts build fresh
# runs instead of mutating (and merging into) the existing neighboring run.
$target = $d.Range($startIdx, $endIdx)
$target.Delete()

# Re-insert the text as three separate runs, split around the missing "t"
# that fixes "seting" -> "setting":
#   1) "...academic se"
#   2) "t"
#   3) "ting. "
$part1 = "with many different tools and software such as Python, Solidworks, and computer vision in a manner and depth that are unavailable in an academic se"
$part2 = "t"
$part3 = "ting. "

$insPoint = $d.Range($startIdx, $startIdx)
$insPoint.InsertAfter($part1)

$full = $d.Content.Text
$p2 = $full.IndexOf($part1) + $part1.Length
$insPoint2 = $d.Range($p2, $p2)
$insPoint2.InsertAfter($part2)

$full = $d.Content.Text
$p3 = $full.IndexOf($part1 + $part2) + ($part1 + $part2).Length
$insPoint3 = $d.Range($p3, $p3)
$insPoint3.InsertAfter($part3)
